$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 134; this shifts rows 134..247 down to 135..248,
# matching the dimension change from A1:R247 to A1:R248.
$ws.Rows("134:134").Insert()

# Populate the newly inserted row 134 with the new record.
$ws.Range("A134").Value = 7
$ws.Range("B134").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C134").Value = "Ñuble"
$ws.Range("D134").Value = 44651
$ws.Range("D134").NumberFormat = $ws.Range("D135").NumberFormat
$ws.Range("E134").Value = 16
$ws.Range("F134").Value = 100112023
$ws.Range("G134").Value = "Brócoli"
$ws.Range("H134").Value = "Sin especificar"
$ws.Range("I134").Value = "Segunda"
$ws.Range("J134").Value = 100
$ws.Range("K134").Value = 650
$ws.Range("L134").Value = 650
$ws.Range("M134").Value = 650
$ws.Range("N134").Value = "$/unidad"
$ws.Range("O134").Value = "Región del Maule"
$ws.Range("P134").Value = 650
$ws.Range("Q134").Value = 1
$ws.Range("R134").Value = "Hortaliza"
